$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: "Curso (semestre ideal): EA (7)" -> "... EA (8)"
# --------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Curso (semestre ideal): EA (7)", $false, $false, $false, $false,
    $false, $true, 1, $false, "Curso (semestre ideal): EA (8)", 2)

# --------------------------------------------------------------------
# Change 2: rework the "Requisitos" bullet list.
#   LOB1207 - Poluição Ambiental I        (Requisito)
#   LOB1208 - Química Analítica Ambiental I (Requisito)   <- removed
#   LOB1255 - Hidrologia Aplicada         (Requisito)
# becomes
#   LOB1212 - Química Analítica Ambiental II (Requisito fraco)
#   LOB1258 - Hidráulica Aplicada            (Requisito fraco)
# --------------------------------------------------------------------

# Locate the "Requisitos" heading paragraph, then the bullet-list
# paragraph that immediately follows it (robust to absolute index).
$reqHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("Requisitos") -and -not $t.Contains("LOB")) {
        $reqHeadingIndex = $i
        break
    }
}
$listPara = $d.Paragraphs.Item($reqHeadingIndex + 1)
$listRange = $listPara.Range

# Drop the trailing paragraph mark from the range so InsertXML only
# touches the run content, leaving the paragraph (and its ListBullet
# style) intact.
$target = $d.Range($listRange.Start, $listRange.End - 1)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOB1212 -  Química Analítica Ambiental II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1258 -  Hidráulica Aplicada  (Requisito fraco)</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag)
